# Generate Report for Handoff
# Update Priority ("low" -> "ht") and Latest Handoff Datetime for the rows that
# were still pending ("Ready for handoff") on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 had Priority "low"; Latest Handoff Datetime moves forward
# from 2016-08-19 20:36:55 to 2016-08-19 20:37:15.
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-19 20:37:15"
}

# de-de sheet: rows 4-7 had Priority "low"; Latest Handoff Datetime moves forward
# from 2016-08-19 20:36:59 to 2016-08-19 20:37:19.
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-19 20:37:19"
}
